$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the TF (H) column scores for each student
$ws.Range("H2").Value = 2.5
$ws.Range("H3").Value = 2.5
$ws.Range("H4").Value = 2
$ws.Range("H5").Value = 2.5
$ws.Range("H6").Value = 2.5
$ws.Range("H6").NumberFormat = "0.00"

# Fill in the Conceito (J) column with "A" for each student
$ws.Range("J2").Value = "A"
$ws.Range("J3").Value = "A"
$ws.Range("J4").Value = "A"
$ws.Range("J5").Value = "A"
$ws.Range("J6").Value = "A"

# Update selection to I2 as in the final state
$ws.Range("I2").Select()
